# NPM Process TC.docx edit script
#
# Per the commit message ("Got rid of the packages folder, which was
# causing conflicts") the two path references to
#   C:\Projects\TypedContract\Code\typed-contract
# are trimmed to
#   C:\Projects\TypedContract\Code\
# and the stray "_GoBack" bookmark that used to sit at the end of the
# "Use the credentials..." bullet is removed from there and re-created
# at the spot of the last real text edit: right before the word
# "directory" in the "Go to the ... directory" bullet.
#
# Implementation notes (quirks of this COM-interop runtime worked
# around below):
#
#  1) Any Range-based text mutation (Delete/Text=/InsertBefore/...)
#     that lands *inside* an existing run causes this engine to
#     coalesce that run with same-formatted neighboring runs in the
#     paragraph, even ones the edit never touched. Dropping a
#     temporary bookmark at a position acts as a hard boundary that
#     blocks that coalescing pass; removing the temporary bookmark
#     again afterwards does not undo the split. We use this
#     "barrier bookmark" trick everywhere we need a run boundary to
#     survive an edit.
#
#  2) A run created by *splitting* pre-existing "preserve" text (e.g.
#     splitting " directory" into " " and "directory" in place) keeps
#     xml:space="preserve" on both halves even when a half no longer
#     needs it. Text that is instead freshly inserted via
#     Range.InsertAfter does not get that attribute. So the final
#     "directory" run is produced by deleting the old word and
#     re-inserting it fresh, rather than just splitting around it.
#
#  3) Adding a bookmark at a zero-length Range sitting at exactly
#     "paragraph end minus one" (the last real character before the
#     paragraph mark) lands the bookmark at document position 0
#     instead of where asked. We avoid ever adding a bookmark while
#     its target position is that last-character-before-mark spot by
#     inserting the new "_GoBack" bookmark *before* deleting the old
#     "directory" word (so the bookmark's insertion point still has
#     real run content after it, i.e. is not at the paragraph's last
#     position yet).

$d = $word.ActiveDocument

function Insert-Barrier([int]$pos, [string]$name) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $r) | Out-Null
}

function Remove-Barrier([string]$name) {
    if ($d.Bookmarks.Exists($name)) {
        $d.Bookmarks($name).Delete()
    }
}

$prefix = "C:\Projects\TypedContract\Code\"
$needle = $prefix + "typed-contract"

# ---------------------------------------------------------------
# 1) Trim "...\Code\typed-contract" -> "...\Code\" in the
#    "Make sure the version number in ..." bullet (first occurrence
#    in the document).
# ---------------------------------------------------------------
$full = $d.Content.Text
$occ1 = $full.IndexOf($needle)
if ($occ1 -lt 0) { throw "first occurrence of path not found" }

Insert-Barrier $occ1 "ZZZ_B1S"
Insert-Barrier ($occ1 + $needle.Length) "ZZZ_B1E"

$suffixStart = $occ1 + $prefix.Length
$suffixEnd = $occ1 + $needle.Length
$d.Range($suffixStart, $suffixEnd).Delete() | Out-Null

Remove-Barrier "ZZZ_B1S"
Remove-Barrier "ZZZ_B1E"

# ---------------------------------------------------------------
# 2) Trim "...\Code\typed-contract" -> "...\Code\" in the
#    "Go to the ... directory" bullet (the remaining occurrence).
# ---------------------------------------------------------------
$full = $d.Content.Text
$occ2 = $full.IndexOf($needle)
if ($occ2 -lt 0) { throw "second occurrence of path not found" }

Insert-Barrier $occ2 "ZZZ_B2S"
Insert-Barrier ($occ2 + $needle.Length) "ZZZ_B2E"

$suffixStart = $occ2 + $prefix.Length
$suffixEnd = $occ2 + $needle.Length
$d.Range($suffixStart, $suffixEnd).Delete() | Out-Null

Remove-Barrier "ZZZ_B2S"
Remove-Barrier "ZZZ_B2E"

# ---------------------------------------------------------------
# 3) Remove the old "_GoBack" bookmark from the "Use the
#    credentials..." bullet.
# ---------------------------------------------------------------
Remove-Barrier "_GoBack"

# ---------------------------------------------------------------
# 4) In the "Go to the ... directory" bullet, split the trailing
#    " directory" text into a " " run followed by a "_GoBack"
#    bookmark followed by a fresh "directory" run.
# ---------------------------------------------------------------
$full = $d.Content.Text
$dirNeedle = $prefix + " directory"
$occ3 = $full.IndexOf($dirNeedle)
if ($occ3 -lt 0) { throw "'" + $dirNeedle + "' not found" }

$spacePos = $occ3 + $prefix.Length   # position of the single space char
$wordPos = $spacePos + 1             # "directory" begins right after it
$wordEnd = $wordPos + "directory".Length

# Barrier in front of the space so trimming the path above doesn't
# leave this area mergeable with the prefix once we touch it below.
Insert-Barrier $spacePos "ZZZ_B3S"

# Drop the new "_GoBack" bookmark *before* removing the old
# "directory" word, while real content still follows the insertion
# point (dodges the "last position in paragraph" bookmark bug).
$bmRange = $d.Range($wordPos, $wordPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Delete the stale "directory" word that now sits right after the
# bookmark.
$d.Range($wordPos, $wordEnd).Delete() | Out-Null

Remove-Barrier "ZZZ_B3S"

# Re-type "directory" fresh so it doesn't inherit xml:space="preserve"
# from the run it used to share with the leading space.
$d.Range($wordPos, $wordPos).InsertAfter("directory") | Out-Null
